$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.161197760840216
$ws.Range("C2").Value = 0.01163354830652435

$ws.Range("B3").Value = 0.1439540370011767
$ws.Range("C3").Value = 0.001672567678019396

$ws.Range("B4").Value = 0.3038949153813413
$ws.Range("C4").Value = 0.01535662063916811

$ws.Range("B5").Value = 0.2703798058130104
$ws.Range("C5").Value = 0.003451387287577847

$ws.Range("B6").Value = 0.2845714814316892
$ws.Range("C6").Value = 0.001870071720553954
